$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Selection change (sheetView) ---
[void]$ws.Range("G14").Select()

# --- Key table (K/L) & NPC "Type" column (I) updates ---
# Type column for Sheep/Duck/Cow rows moves from type 3 to type 4
$ws.Range("I4").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("I6").Value = 4

# --- Map Construction table (N/O/P) ---
# Row 11 becomes "Map Doors" (new concept), min defaults to 0
# (write this shared string before "Is a door." so new sharedStrings entries
# land in the same order as the target: 68="Map Doors", 69="Is a door.")
$ws.Range("N11").Value = "Map Doors"
$ws.Range("O11").Value = 0

# Key table: type 3 is now "Is a door." (was "Is an NPC.")
$ws.Range("L7").Value = "Is a door."

# Key table: new type 4 row = "Is an NPC."
$ws.Range("K8").Value = 4
$ws.Range("K8").HorizontalAlignment = -4131
$ws.Range("L8").Value = "Is an NPC."

# Row 12 (new) becomes what used to be row 11: "Map data"
$ws.Range("N12").Value = "Map data"
$ws.Range("O12").Value = 1
$ws.Range("P12").Value = 3600

# Row 13's note moves down to row 14; row 13 no longer holds it
$ws.Range("N13").ClearContents()
$ws.Range("N14").Value = "Note: See example below on how the document MUST be formatted."

# Entrance ("X") Type changes from 0 to 3 (now represents a door)
$ws.Range("D14").Value = 3

# --- Picture anchor/position update ---
$shp = $ws.Shapes.Item(1)
$shp.Left = 741.75
$shp.Top = 225
$shp.Width = 216.75
$shp.Height = 201
